$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Extend the existing "What to include" help text (HelpBoxOptions!A3) with
#    a "Do not include" section.
# ---------------------------------------------------------------------------
$helpBox = $wb.Worksheets.Item("HelpBoxOptions")

$existingText = $helpBox.Range("A3").Value()
$newText = $existingText + "`nDo not include:`nCurrent rent or mortgage payments that you" + [char]0x2019 + "ll no longer pay after buying your new home`nYour new monthly mortgage payment`nCredit card balances that you pay off in full each month`nMonthly utilities, groceries, and other costs of living"
$helpBox.Range("A3").Value = $newText
$helpBox.Range("A3").RowHeight = 182.25

$helpBox.Range("A5").Value = "Property taxes and interest rates can vary by location. Enter the ZIP code where you are looking for homes and we will automatically add property taxes and interest rates for your area."
$helpBox.Range("A4").Value = "Your credit score affects your loan eligibility and your interest rate. Typically, the higher your credit score, the lower your interest rate."

$helpBox.Range("A6").Select()

# ---------------------------------------------------------------------------
# 2. MortgageQuestions is no longer the tab that is selected when the file
#    is opened.
# ---------------------------------------------------------------------------
$mortgageQuestions = $wb.Worksheets.Item("MortgageQuestions")
$mortgageQuestions.Select()

# ---------------------------------------------------------------------------
# 3. Add a new "MortgageDropdown" worksheet at the end of the workbook with
#    the new dropdown copy.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "MortgageDropdown"

$newSheet.Range("A1").Value = "Buy"
$newSheet.Range("A2").Value = "Rent"
$newSheet.Range("A3").Value = "Mortgage"
$newSheet.Range("A4").Value = "Saved Homes"
$newSheet.Range("A5").Value = "Saved Searches"
$newSheet.Range("A6").Value = "Sign up or Log in"

$newSheet.Activate()
$newSheet.Range("A1:C6").Select()
